# Replace each three-digit x one-digit multiplication problem in the practice
# table with the updated problem from the same cell position, per the commit diff.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$cell = $t.Cell(1, 1)
$cell.Range.Text = "435×3="

$cell = $t.Cell(1, 2)
$cell.Range.Text = "992×4="

$cell = $t.Cell(1, 3)
$cell.Range.Text = "839×4="

$cell = $t.Cell(1, 4)
$cell.Range.Text = "916×2="

$cell = $t.Cell(1, 5)
$cell.Range.Text = "329×2="

$cell = $t.Cell(5, 1)
$cell.Range.Text = "157×4="

$cell = $t.Cell(5, 2)
$cell.Range.Text = "637×4="

$cell = $t.Cell(5, 3)
$cell.Range.Text = "873×8="

$cell = $t.Cell(5, 4)
$cell.Range.Text = "423×7="

$cell = $t.Cell(5, 5)
$cell.Range.Text = "173×6="

$cell = $t.Cell(10, 1)
$cell.Range.Text = "261×2="

$cell = $t.Cell(10, 2)
$cell.Range.Text = "774×3="

$cell = $t.Cell(10, 3)
$cell.Range.Text = "403×9="

$cell = $t.Cell(10, 4)
$cell.Range.Text = "630×7="

$cell = $t.Cell(10, 5)
$cell.Range.Text = "298×7="

$cell = $t.Cell(15, 1)
$cell.Range.Text = "997×7="

$cell = $t.Cell(15, 2)
$cell.Range.Text = "199×3="

$cell = $t.Cell(15, 3)
$cell.Range.Text = "821×2="

$cell = $t.Cell(15, 4)
$cell.Range.Text = "276×6="

$cell = $t.Cell(15, 5)
$cell.Range.Text = "414×9="

$cell = $t.Cell(20, 1)
$cell.Range.Text = "528×6="

$cell = $t.Cell(20, 2)
$cell.Range.Text = "361×6="

$cell = $t.Cell(20, 3)
$cell.Range.Text = "412×2="

$cell = $t.Cell(20, 4)
$cell.Range.Text = "911×7="

$cell = $t.Cell(20, 5)
$cell.Range.Text = "609×5="
